$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$features = @(
    @('YB/LU', 1572),
    @('TM/YB', 2218.200000000001),
    @('ER/TM', 2085),
    @('TB/DY', 2072),
    @('CE**EU*EU', 1854.4),
    @('Y/LA/TB', 1840.799999999999),
    @('DY/HO', 1790),
    @('TH/U', 1758),
    @('EU/GD', 1698.2),
    @('Y/CE*/TH', 1642.6),
    @('YB/LU/HF', 1621.4),
    @('TA/U', 1615.6),
    @('HO-LU', 1606),
    @('NB/LA/CE*', 1581.2),
    @('ER/YB', 1577.2),
    @('PR/ND/EU*', 1535.599999999999),
    @('LU/TH', 1517.2),
    @('EU**EU**TA', 1497),
    @('Y/CE*/EU', 1451.2),
    @('TI*Y*EU*', 1444.6),
    @('NB/LU/TA', 1426.2),
    @('GD/TB', 1423.4),
    @('Y/YB/TH', 1421.4),
    @('PR/ND/TA', 1396.2),
    @('Y/CE/EU*', 1390.2),
    @('CE/HF/TH', 1389.4),
    @('HF-U', 1389),
    @('TH-U', 1377.6),
    @('CE/U', 1360.8),
    @('CE/EU*/EU*', 1348),
    @('CE/EU*/TH', 1345.2),
    @('Y/DY/HF', 1345),
    @('CE-ND-TA', 1341.4),
    @('EU**HF*TA', 1332),
    @('YB/LU/TA', 1318),
    @('CE/DY/U', 1310.2),
    @('PR/EU/TA', 1304.2),
    @('TI/NB/EU*', 1302),
    @('TI*EU**U', 1299.2),
    @('TI*TI*U', 1297.8),
    @('ER/TM/HF', 1294),
    @('CE/GD', 1287),
    @('TI/HF/HF', 1285.8),
    @('HF-TH-U', 1280.8),
    @('NB/CE/CE', 1271.2),
    @('CE/TA/U', 1268.000000000001),
    @('TI/EU*/HF', 1259.2),
    @('Y/CE/U', 1252),
    @('TI*EU*HF', 1247.8),
    @('TI*Y*TA', 1239),
    @('SM/GD', 1220.8),
    @('EU**HF*HF', 1217),
    @('TI*TI*CE', 1211),
    @('Y*CE*EU*', 1177.2),
    @('LA+EU*+TA', 1171.6),
    @('ND-EU-EU*', 1161.8),
    @('CE/EU*', 1157.6),
    @('ER-LU-LU', 1157),
    @('TI*NB*TA', 1154.8),
    @('ER/YB/HF', 1152.2),
    @('DY-LU-LU', 1144.6),
    @('Y/ER', 1133.2),
    @('CE/GD/HF', 1114.2),
    @('SM/EU/TA', 1105.4),
    @('CE*EU**HF', 1102),
    @('Y/HO', 1100.8),
    @('HO-LU-TA', 1096.2),
    @('Y/NB', 1091.2),
    @('ER/LU', 1090.6),
    @('TI*TI*TA', 1085.4),
    @('CE**ND*HF', 1081.8),
    @('Y/NB/ER', 1078.6),
    @('EU/GD/HF', 1076.6),
    @('DY/LU', 1074.4),
    @('CE/EU/TH', 1071.4),
    @('TI*HF*TA', 1066.2),
    @('CE/DY/DY', 1061.8),
    @('HO/TM', 1053.4),
    @('HO/ER', 1048.799999999999),
    @('TM/YB/HF', 1047.8),
    @('TI*HF*HF', 1042),
    @('NB/HF/TH', 1031.8),
    @('TM/LU', 1028),
    @('Y/TM', 1027.6),
    @('TI/EU*/TA', 1025.4),
    @('NB*EU**U', 1023.6),
    @('CE/EU/EU*', 1011.6),
    @('CE/TH', 1010.4),
    @('CE/EU/TA', 1001.8),
    @('NB/HF/U', 999.6),
    @('EU*/TA/TA', 992.2),
    @('CE/TA/TH', 985.1999999999999),
    @('CE/HF/HF', 975.2),
    @('Y/EU*/U', 972.6),
    @('EU/EU*/GD', 969.9999999999999),
    @('TI*TI*DY', 939.7999999999998),
    @('Y/NB/EU*', 936.6000000000003),
    @('Y/SM', 935.2000000000002),
    @('CE/EU*/U', 934.6000000000001),
    @('Y/NB/YB', 931.8000000000001),
    @('Y/TH', 930.7999999999996),
    @('Y*EU**HF', 925.2000000000004),
    @('CE/SM/HF', 924.6),
    @('DY-LU-TA', 923),
    @('ND-EU-EU', 905.7999999999998),
    @('SM/GD/HF', 898.6),
    @('Y/NB/CE', 897.8000000000002),
    @('NB/TH/U', 892.9999999999997),
    @('Y*HF*TA', 884.1999999999999),
    @('Y/YB', 883.1999999999998),
    @('TB/DY/HF', 849.9999999999998),
    @('TI/HF', 843.6000000000001),
    @('Y/DY', 838),
    @('Y/DY/TA', 829.8000000000002),
    @('CE/EU/HF', 820.9999999999999),
    @('NB-TA-TA', 772.3999999999997),
    @('NB/HF/TA', 702.8),
    @('TA*U', 633.2000000000002),
    @('TB/LU/HF', 571.3999999999999),
    @('NB/TA', 511.9999999999999)
)

for ($i = 0; $i -lt $features.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $features[$i][0]
    $ws.Cells.Item($row, 2).Value = $features[$i][1]
}

Write-Host "Done updating feature importance rows."